$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 7.293566666666666
$ws.Range("H2").Value = 21.8807
$ws.Range("I2").Value = 0.2546895690137356
$ws.Range("J2").Value = 0.2546895690137356
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 50.45577166666666
$ws.Range("N2").Value = 151.367315
$ws.Range("O2").Value = 0.18996722124408
$ws.Range("P2").Value = 0.1899672212440799
$ws.Range("Q2").Value = 368.0025343689444
$ws.Range("R2").Value = 3312.0228093205
$ws.Range("S2").Value = 0.04838266970539169
$ws.Range("T2").Value = 0.04838266970539167

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 7.293566666666666
$ws.Range("H3").Value = 21.8807
$ws.Range("I3").Value = 0.2546895690137356
$ws.Range("J3").Value = 0.2546895690137356
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 149.2656146666667
$ws.Range("N3").Value = 447.796844
$ws.Range("O3").Value = 0.5619887102876124
$ws.Range("P3").Value = 0.5619887102876123
$ws.Range("Q3").Value = 1088.678711612311
$ws.Range("R3").Value = 9798.1084045108
$ws.Range("S3").Value = 0.1431326624137371
$ws.Range("T3").Value = 0.1431326624137371

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 7.293566666666666
$ws.Range("H4").Value = 21.8807
$ws.Range("I4").Value = 0.2546895690137356
$ws.Range("J4").Value = 0.2546895690137356
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 13.28522433333333
$ws.Range("N4").Value = 39.855673
$ws.Range("O4").Value = 0.05001919635439596
$ws.Range("P4").Value = 0.05001919635439595
$ws.Range("Q4").Value = 96.89666935678888
$ws.Range("R4").Value = 872.0700242110998
$ws.Range("S4").Value = 0.01273936756191452
$ws.Range("T4").Value = 0.01273936756191452

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 7.293566666666666
$ws.Range("H5").Value = 21.8807
$ws.Range("I5").Value = 0.2546895690137356
$ws.Range("J5").Value = 0.2546895690137356
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 17.42465066666666
$ws.Range("N5").Value = 52.27395199999999
$ws.Range("O5").Value = 0.06560423830525379
$ws.Range("P5").Value = 0.06560423830525378
$ws.Range("Q5").Value = 127.0878512807111
$ws.Range("R5").Value = 1143.7906615264
$ws.Range("S5").Value = 0.01670871517943949
$ws.Range("T5").Value = 0.01670871517943949

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 7.293566666666666
$ws.Range("H6").Value = 21.8807
$ws.Range("I6").Value = 0.2546895690137356
$ws.Range("J6").Value = 0.2546895690137356
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 3.820091
$ws.Range("N6").Value = 11.460273
$ws.Range("O6").Value = 0.01438273656706242
$ws.Range("P6").Value = 0.01438273656706242
$ws.Range("Q6").Value = 27.86208838123333
$ws.Range("R6").Value = 250.7587954311
$ws.Range("S6").Value = 0.003663132977503223
$ws.Range("T6").Value = 0.003663132977503223

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 7.293566666666666
$ws.Range("H7").Value = 21.8807
$ws.Range("I7").Value = 0.2546895690137356
$ws.Range("J7").Value = 0.2546895690137356
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 31.35116233333333
$ws.Range("N7").Value = 94.053487
$ws.Range("O7").Value = 0.1180378972415954
$ws.Range("P7").Value = 0.1180378972415954
$ws.Range("Q7").Value = 228.6617925556556
$ws.Range("R7").Value = 2057.9561330009
$ws.Range("S7").Value = 0.03006302117574954
$ws.Range("T7").Value = 0.03006302117574954

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 9.684806
$ws.Range("H8").Value = 29.054418
$ws.Range("I8").Value = 0.3381910632824783
$ws.Range("J8").Value = 0.3381910632824783
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 50.45577166666666
$ws.Range("N8").Value = 151.367315
$ws.Range("O8").Value = 0.18996722124408
$ws.Range("P8").Value = 0.1899672212440799
$ws.Range("Q8").Value = 488.6543601719633
$ws.Range("R8").Value = 4397.88924154767
$ws.Range("S8").Value = 0.0642452165413532
$ws.Range("T8").Value = 0.06424521654135319

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 9.684806
$ws.Range("H9").Value = 29.054418
$ws.Range("I9").Value = 0.3381910632824783
$ws.Range("J9").Value = 0.3381910632824783
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 149.2656146666667
$ws.Range("N9").Value = 447.796844
$ws.Range("O9").Value = 0.5619887102876124
$ws.Range("P9").Value = 0.5619887102876123
$ws.Range("Q9").Value = 1445.608520517421
$ws.Range("R9").Value = 13010.47668465679
$ws.Range("S9").Value = 0.1900595594849163
$ws.Range("T9").Value = 0.1900595594849163

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 9.684806
$ws.Range("H10").Value = 29.054418
$ws.Range("I10").Value = 0.3381910632824783
$ws.Range("J10").Value = 0.3381910632824783
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 13.28522433333333
$ws.Range("N10").Value = 39.855673
$ws.Range("O10").Value = 0.05001919635439596
$ws.Range("P10").Value = 0.05001919635439595
$ws.Range("Q10").Value = 128.6648203348126
$ws.Range("R10").Value = 1157.983383013314
$ws.Range("S10").Value = 0.01691604519962823
$ws.Range("T10").Value = 0.01691604519962823

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 9.684806
$ws.Range("H11").Value = 29.054418
$ws.Range("I11").Value = 0.3381910632824783
$ws.Range("J11").Value = 0.3381910632824783
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 17.42465066666666
$ws.Range("N11").Value = 52.27395199999999
$ws.Range("O11").Value = 0.06560423830525379
$ws.Range("P11").Value = 0.06560423830525378
$ws.Range("Q11").Value = 168.7543613244373
$ws.Range("R11").Value = 1518.789251919936
$ws.Range("S11").Value = 0.02218676710829087
$ws.Range("T11").Value = 0.02218676710829087

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 9.684806
$ws.Range("H12").Value = 29.054418
$ws.Range("I12").Value = 0.3381910632824783
$ws.Range("J12").Value = 0.3381910632824783
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 3.820091
$ws.Range("N12").Value = 11.460273
$ws.Range("O12").Value = 0.01438273656706242
$ws.Range("P12").Value = 0.01438273656706242
$ws.Range("Q12").Value = 36.996840237346
$ws.Range("R12").Value = 332.9715621361141
$ws.Range("S12").Value = 0.004864112972526622
$ws.Range("T12").Value = 0.004864112972526622

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 9.684806
$ws.Range("H13").Value = 29.054418
$ws.Range("I13").Value = 0.3381910632824783
$ws.Range("J13").Value = 0.3381910632824783
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 31.35116233333333
$ws.Range("N13").Value = 94.053487
$ws.Range("O13").Value = 0.1180378972415954
$ws.Range("P13").Value = 0.1180378972415954
$ws.Range("Q13").Value = 303.6299250728407
$ws.Range("R13").Value = 2732.669325655566
$ws.Range("S13").Value = 0.03991936197576306
$ws.Range("T13").Value = 0.03991936197576305

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 11.65871166666667
$ws.Range("H14").Value = 34.976135
$ws.Range("I14").Value = 0.4071193677037862
$ws.Range("J14").Value = 0.4071193677037861
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 50.45577166666666
$ws.Range("N14").Value = 151.367315
$ws.Range("O14").Value = 0.18996722124408
$ws.Range("P14").Value = 0.1899672212440799
$ws.Range("Q14").Value = 588.2492937808361
$ws.Range("R14").Value = 5294.243644027524
$ws.Range("S14").Value = 0.0773393349973351
$ws.Range("T14").Value = 0.07733933499733507

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 11.65871166666667
$ws.Range("H15").Value = 34.976135
$ws.Range("I15").Value = 0.4071193677037862
$ws.Range("J15").Value = 0.4071193677037861
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 149.2656146666667
$ws.Range("N15").Value = 447.796844
$ws.Range("O15").Value = 0.5619887102876124
$ws.Range("P15").Value = 0.5619887102876123
$ws.Range("Q15").Value = 1740.244763146438
$ws.Range("R15").Value = 15662.20286831794
$ws.Range("S15").Value = 0.228796488388959
$ws.Range("T15").Value = 0.228796488388959

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 11.65871166666667
$ws.Range("H16").Value = 34.976135
$ws.Range("I16").Value = 0.4071193677037862
$ws.Range("J16").Value = 0.4071193677037861
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 13.28522433333333
$ws.Range("N16").Value = 39.855673
$ws.Range("O16").Value = 0.05001919635439596
$ws.Range("P16").Value = 0.05001919635439595
$ws.Range("Q16").Value = 154.8885999293172
$ws.Range("R16").Value = 1393.997399363855
$ws.Range("S16").Value = 0.02036378359285321
$ws.Range("T16").Value = 0.0203637835928532

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 11.65871166666667
$ws.Range("H17").Value = 34.976135
$ws.Range("I17").Value = 0.4071193677037862
$ws.Range("J17").Value = 0.4071193677037861
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 17.42465066666666
$ws.Range("N17").Value = 52.27395199999999
$ws.Range("O17").Value = 0.06560423830525379
$ws.Range("P17").Value = 0.06560423830525378
$ws.Range("Q17").Value = 203.1489780150578
$ws.Range("R17").Value = 1828.34080213552
$ws.Range("S17").Value = 0.02670875601752343
$ws.Range("T17").Value = 0.02670875601752342

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 11.65871166666667
$ws.Range("H18").Value = 34.976135
$ws.Range("I18").Value = 0.4071193677037862
$ws.Range("J18").Value = 0.4071193677037861
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 3.820091
$ws.Range("N18").Value = 11.460273
$ws.Range("O18").Value = 0.01438273656706242
$ws.Range("P18").Value = 0.01438273656706242
$ws.Range("Q18").Value = 44.53733950942834
$ws.Range("R18").Value = 400.836055584855
$ws.Range("S18").Value = 0.005855490617032577
$ws.Range("T18").Value = 0.005855490617032576

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 11.65871166666667
$ws.Range("H19").Value = 34.976135
$ws.Range("I19").Value = 0.4071193677037862
$ws.Range("J19").Value = 0.4071193677037861
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 31.35116233333333
$ws.Range("N19").Value = 94.053487
$ws.Range("O19").Value = 0.1180378972415954
$ws.Range("P19").Value = 0.1180378972415954
$ws.Range("Q19").Value = 365.5141620591939
$ws.Range("R19").Value = 3289.627458532745
$ws.Range("S19").Value = 0.0480555140900828
$ws.Range("T19").Value = 0.04805551409008279
